# إضافة صف جديد في Card9
#
# 1) Row 29: cells B29:K29 were empty placeholders -> fill them with the
#    literal text "nan" (matching the sheet's existing convention for
#    "no data" elsewhere in the table).
# 2) Row 30: append the newest service-log entry.
#
# Numeric-looking ("9", "1001", ...) and date-looking ("9/12/2025") text
# needs to be entered as plain TEXT (this sheet stores every column as
# text, never as a real number/date). Typing those strings straight into
# a General-formatted cell would make Excel auto-convert them to a real
# number / date serial, so each such cell is briefly switched to the
# Text number format before the value is typed in, then the format is
# cleared again so no extra cell styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

function Set-BlankCell($row, $col) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.ClearFormats()
}

# --- Row 29: cells B29:K29 were empty inline strings; set them to "nan" ---
Set-TextValue 29 2  "nan"   # B29
Set-TextValue 29 3  "nan"   # C29
Set-TextValue 29 4  "nan"   # D29
Set-TextValue 29 5  "nan"   # E29
Set-TextValue 29 6  "nan"   # F29
Set-TextValue 29 7  "nan"   # G29
Set-TextValue 29 8  "nan"   # H29
Set-TextValue 29 9  "nan"   # I29
Set-TextValue 29 10 "nan"   # J29
Set-TextValue 29 11 "nan"   # K29

# --- Row 30: new row of data ---
Set-TextValue 30 1  "9"                 # A30 - card
Set-TextValue 30 2  "1001"              # B30 - Min_Tones
Set-TextValue 30 3  "1150"              # C30 - Max_Tones
Set-TextValue 30 4  "1035"              # D30 - Tones
Set-BlankCell 30 5                      # E30 - Revolving flats(x)
Set-BlankCell 30 6                      # F30 - 1.carding elemnt(o)
Set-BlankCell 30 7                      # G30 - licker_in carding element(o)
$ws.Cells.Item(30, 8).Value = "✅"       # H30 - Doffer carding element(o)
Set-BlankCell 30 9                      # I30 - cylinder(X)
Set-BlankCell 30 10                     # J30 - doffer(X)
Set-BlankCell 30 11                     # K30 - Revolving flats(o)
Set-TextValue 30 12 "9/12/2025"         # L30 - Date
Set-BlankCell 30 13                     # M30 - Event
Set-BlankCell 30 14                     # N30 - Correction
$ws.Cells.Item(30, 15).Value = "م محمد عبدالله "   # O30 - Serviced by (trailing space preserved)

$wb.Save()
